# Apply attendance classification updates to Sheet1.
# Columns: D=Total Attendance Count, E=Real, F=Duplicate, G=Invalid, H=Absent
# For each date row, set the appropriate status column(s) from 0 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 1

$ws.Range("H7").Value = 1

$ws.Range("H8").Value = 1

$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 1

$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1

$ws.Range("H11").Value = 1

$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1

$ws.Range("H13").Value = 1

$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 1

$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 1

$ws.Range("H16").Value = 1

$ws.Range("H17").Value = 1

$ws.Range("H18").Value = 1
